# Append a new row (row 97) of data to each of the four worksheets,
# mirroring the layout/formatting of the existing last row (row 96).

$wb = $excel.ActiveWorkbook

$rowData = @{
    1 = @{
        A = 45883.49545138889
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x00"
        E = "0xf"
        F = 380
        G = "7.598631275147109e+23"
        H = 256
        I = 15
    }
    2 = @{
        A = 45883.49545138889
        B = "0x01,0x90"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x14"
        E = "0xe"
        F = 400
        G = "5.68432987514711e+23"
        H = 276
        I = 14
    }
    3 = @{
        A = 45883.49545138889
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x5C"
        E = "0x3"
        F = 110
        G = "5.68631262647114e+23"
        H = 92
        I = 3
    }
    4 = @{
        A = 45883.49545138889
        B = "0x00,0x6e"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x5A"
        E = "0x3"
        F = 110
        G = "9.85046333984776e+23"
        H = 90
        I = 3
    }
}

foreach ($sheetIndex in 1..4) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    $data = $rowData[$sheetIndex]

    $newRow = 97
    $sourceRow = 96

    $ws.Cells.Item($newRow, 1).Value = $data.A
    $ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($sourceRow, 1).NumberFormat

    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E

    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = [double]$data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}
